$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9083
$ws.Range("C3:C8").Value = 8793
$ws.Range("C9:C16").Value = 8777
$ws.Range("C17:C20").Value = 8183
$ws.Range("C21:C45").Value = 7861
$ws.Range("C46:C49").Value = 7672
$ws.Range("C50:C252").Value = 7622

Write-Host "Done updating Fitness column C2:C252"
